$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 4 and 5 entirely - they are no longer part of the table.
$ws.Rows("4:5").Delete()

# Row 2: SJ2 / param / das@dfg.com / 413548 (as a real number)
$ws.Range("A2").Value = "SJ2"
$ws.Range("B2").Value = "param"
$ws.Range("C2").Value = "das@dfg.com"
$ws.Range("D2").Value = 413548

# Row 3: SJ3 / alok / alohg / 65674537 (kept as text, not a number)
$ws.Range("A3").Value = "SJ3"
$ws.Range("B3").Value = "alok"
$ws.Range("C3").Value = "alohg"

# D3 must stay a text value ("65674537") rather than being auto-converted to a
# number. Build it in a scratch cell formatted as Text, copy/paste-special
# just the value into D3 (so D3 itself keeps the default, unstyled format),
# then clean the scratch cell back up.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "65674537"
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$scratch.Clear()
